$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(34).Insert()
$ws.Cells.Item(34, 1).Value = "electrocoagulation"
$ws.Cells.Item(34, 2).Value = "electrocoagulation_zo"
$ws.Cells.Item(34, 3).Value = "ElectrocoagulationZO"
$ws.Cells.Item(34, 4).Value = "non-basic"
$ws.Cells.Item(34, 5).Value = "SIDO"
$ws.Cells.Item(34, 6).Value = $false
$ws.Cells.Item(34, 7).Value = "f(x)"
$ws.Cells.Item(34, 8).Value = "cost_electrocoagulation"
$ws.Cells.Item(34, 9).Formula = '=IF(E34="SIDO","single-input, double-output",IF(E34="SISO","single-input, single-output",IF(E34="PT","pass-through",IF(E34="DISO","double-input, single-output",IF(E34="SIDO reactive","reactive single-inlet, double-outlet","")))))'
$ws.Cells.Item(34, 10).Formula = '=IF(E34="SIDO","sido_methods",IF(E34="SISO","siso_methods",IF(E34="PT","pt_methods",IF(E34="DISO","diso_methods",IF(E34="SIDO reactive","sidor_methods","")))))'
